$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 5 (rows 14-22): Rannsóknir, Friday (G16) updated hours
$ws.Range("G16").Value = 0.667

# Week 6 (rows 24-32): new hours entered
$ws.Range("C27").Value = 0.75
$ws.Range("C31").Value = 0.25

# Update the selected cell/active cell to reflect the saved view state
$ws.Range("G17").Select()
